$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the "Featured Video" card in B7 with a new "featured_blog" card
# written by justaashir.
$newCard = @"
type: featured_blog
width: 2
height: 1
h3: Rules of being a good desi
p: Some rules to follow if you want to lit Pakistan brighter. We here at zakatlists are bounded by these rules. 😀
date: 6 Apr 2020
author: <a href=https://justaashir.com target=_blank>Aashir</a>
"@

$ws.Range("B7").Value = $newCard

# The edit leaves the active selection on the cell that was just edited.
$ws.Range("B7").Select()
